# Update cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.234.81"
$ws.Range("D3").Value = "2.961.75"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'383.01"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'103.29"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "'36.65"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'0.0841"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "3.426.53"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "'18.06"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "2.942.06"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "'0.991"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").Value = "51.178.90"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'3.22"
$ws.Range("E19").Value = "  -6.01%  "
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").Value = "'12.60"
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'68.50"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'262.42"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").Value = "'2.92"
$ws.Range("E25").Value = "  +3.64%  "
$ws.Range("D26").Value = "'8.41"
$ws.Range("E26").Value = "  +13.11%  "
$ws.Range("D27").Value = "'7.81"
$ws.Range("E27").Value = "  +3.74%  "
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +8.73%  "
$ws.Range("D31").Value = "'25.73"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").Value = "'9.82"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("D34").Value = "'33.97"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'50.41"
$ws.Range("E35").Value = "  -3.81%  "
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "'16.82"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "'2.55"
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").Value = "'0.116"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").Value = "'121.65"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "'21.37"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").Value = "'0.274"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").Value = "'3.25"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").Value = "2.014.01"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "'0.0348"
$ws.Range("E50").Value = "  +6.79%  "
$ws.Range("D51").Value = "'2.11"
$ws.Range("E51").Value = "  +14.13%  "
